# Applies the edits described by the target diff:
#  1. Merge the "Fridlysta arter" intro paragraph + two ListBullet
#     paragraphs ("Skogsfru (NT, §8)" / "Skogsrör (§7)") into a single
#     paragraph.
#  2. Drop the italic run-properties on the first "Kommentar:" follow-up
#     sentence (under 1.3.1) and append a period to its text.
#  3. Trim the trailing space from the "6.4" paragraph's text.
#  4. Remove the two "6.4.1" biotope paragraphs and renumber the
#     remaining "6.4.1" paragraph to "6.4.3".
#  5. Update the title-page header date from 2023-10-22 to 2023-10-25.

$d = $word.ActiveDocument
$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaIndex($r) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($r.Start -ge $p.Range.Start -and $r.Start -lt $p.Range.End) {
            return $i
        }
    }
    return -1
}

# --- 1. Merge "Fridlysta arter" paragraphs ---------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Följande fridlysta arter har sina livsmiljöer och växtplatser i den avverkningsanmälda skogen: "
$found = $find.Execute()
if ($found) {
    $idx = Get-ParaIndex $find.Parent
    $startPara = $d.Paragraphs.Item($idx)
    $endPara = $d.Paragraphs.Item($idx + 2)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $xml = "<w:p $wordNs><w:r><w:t>Följande fridlysta arter har sina livsmiljöer och växtplatser i den avverkningsanmälda skogen: skogsfru (NT, §8) och skogsrör (§7).</w:t></w:r></w:p>"
    $r.InsertXML($xml)
}

# --- 2. First "Kommentar:" sentence - drop italics, add trailing period ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen"
$found = $find.Execute()
if ($found) {
    $idx = Get-ParaIndex $find.Parent
    $p = $d.Paragraphs.Item($idx)
    $xml = "<w:p $wordNs><w:pPr><w:pStyle w:val=`"ListBullet`"/></w:pPr><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t xml:space=`"preserve`">Kommentar: </w:t></w:r><w:r><w:t>I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen.</w:t></w:r></w:p>"
    $p.Range.InsertXML($xml)
}

# --- 3. Trim trailing space on the "6.4" paragraph -------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("ekologiska krav hos sällsynta och hotade arter utanför skogsbruksenhetens gränser när beslut om åtgärder inom skogsbruksenheten ska fattas. ", $true, $false, $false, $false, $false, $true, 1, $false, "ekologiska krav hos sällsynta och hotade arter utanför skogsbruksenhetens gränser när beslut om åtgärder inom skogsbruksenheten ska fattas.", 2)

# --- 4. Remove the two "6.4.1" biotope paragraphs, renumber to 6.4.3 ------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Följande biotoper undantas från alla skogsbruksåtgärder"
$found = $find.Execute()
if ($found) {
    $idx = Get-ParaIndex $find.Parent
    $startPara = $d.Paragraphs.Item($idx)
    $endPara = $d.Paragraphs.Item($idx + 2)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $xml = "<w:p $wordNs><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=`"preserve`">6.4.3 </w:t></w:r><w:r><w:t>Bevarandeåtgärder genomförs för de kända förekomster av rödlistade arter som påverkas av skogsbruk.</w:t></w:r></w:p>"
    $r.InsertXML($xml)
}

# --- 5. Update title-page header date --------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(2)
$xml = "<w:p $wordNs xmlns:w14=`"http://schemas.microsoft.com/office/word/2010/wordml`" w14:paraId=`"042F7434`" w14:textId=`"77777777`" w:rsidR=`"000942A1`" w:rsidRDefault=`"000942A1`"><w:pPr><w:pStyle w:val=`"Header`"/><w:jc w:val=`"left`"/></w:pPr><w:r><w:tab/><w:tab/><w:t>2023-10-25</w:t><w:br/><w:br/></w:r><w:r><w:t>Till:</w:t><w:br/></w:r><w:r><w:t>Kopia: Revisor xx och FSC</w:t><w:br/></w:r></w:p>"
$hdr.Range.InsertXML($xml)

Write-Output "Done."
